# Update the "Fruta, Terminal Hortofrutícola Agro Chillán - Pomelo" sheet.
# The weekly refresh re-shuffled the Fecha/Volumen/Precio rows and appended
# one new observation (row 17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that stay constant across every data row and are reused for the
# newly appended row 17.
$A = 7
$B = "Terminal Hortofrutícola Agro Chillán"
$C = "Ñuble"
$E = 16
$F = "Fruta"
$G = 100102
$H = "Cítricos"
$I = 100102006
$J = "Pomelo"
$K = "Start Ruby"
$L = "Primera"
$R = "Región de O'Higgins"
$T = 14

# Target values (after the edit) for the per-row columns: Fecha (D),
# Volumen (M), Precio minimo (N), Precio maximo (O), Precio promedio (P),
# Unidad de comercializacion (Q) and Precio $/Kg (S).
$rowsData = @(
    @{Row=2;  D=45167; M=50; N=16000; O=16000; P=16000; Q='$/caja 14 kilos empedrada'; S=1143}
    @{Row=3;  D=45138; M=50; N=14000; O=14000; P=14000; Q='$/caja 14 kilos granel';    S=1000}
    @{Row=4;  D=45196; M=30; N=15000; O=15000; P=15000; Q='$/caja 14 kilos empedrada'; S=1071}
    @{Row=5;  D=44210; M=70; N=10000; O=11000; P=10357; Q='$/caja 14 kilos empedrada'; S=740}
    @{Row=6;  D=44216; M=55; N=11000; O=12000; P=11545; Q='$/caja 14 kilos empedrada'; S=825}
    @{Row=7;  D=44229; M=55; N=11000; O=12000; P=11364; Q='$/caja 14 kilos empedrada'; S=812}
    @{Row=8;  D=44172; M=90; N=8500;  O=9000;  P=8806;  Q='$/caja 14 kilos empedrada'; S=629}
    @{Row=9;  D=44253; M=90; N=12000; O=13000; P=12667; Q='$/caja 14 kilos empedrada'; S=905}
    @{Row=10; D=45155; M=60; N=15000; O=15000; P=15000; Q='$/caja 14 kilos empedrada'; S=1071}
    @{Row=11; D=45142; M=30; N=15000; O=15000; P=15000; Q='$/caja 14 kilos empedrada'; S=1071}
    @{Row=12; D=45142; M=30; N=14000; O=14000; P=14000; Q='$/caja 14 kilos granel';    S=1000}
    @{Row=13; D=44181; M=65; N=9000;  O=10000; P=9462;  Q='$/caja 14 kilos empedrada'; S=676}
    @{Row=14; D=45140; M=30; N=15000; O=15000; P=15000; Q='$/caja 14 kilos granel';    S=1071}
    @{Row=15; D=45194; M=60; N=15000; O=15000; P=15000; Q='$/caja 14 kilos granel';    S=1071}
    @{Row=16; D=45152; M=60; N=16000; O=16000; P=16000; Q='$/caja 14 kilos empedrada'; S=1143}
    @{Row=17; D=44232; M=60; N=11000; O=12000; P=11583; Q='$/caja 14 kilos empedrada'; S=827}
)

foreach ($rd in $rowsData) {
    $r = $rd.Row

    # Columns that are constant for every row; row 17 doesn't exist yet, so
    # fill it in completely (harmless no-op writes for the existing rows).
    $ws.Cells.Item($r, 1).Value = $A
    $ws.Cells.Item($r, 2).Value = $B
    $ws.Cells.Item($r, 3).Value = $C
    $ws.Cells.Item($r, 5).Value = $E
    $ws.Cells.Item($r, 6).Value = $F
    $ws.Cells.Item($r, 7).Value = $G
    $ws.Cells.Item($r, 8).Value = $H
    $ws.Cells.Item($r, 9).Value = $I
    $ws.Cells.Item($r, 10).Value = $J
    $ws.Cells.Item($r, 11).Value = $K
    $ws.Cells.Item($r, 12).Value = $L
    $ws.Cells.Item($r, 18).Value = $R
    $ws.Cells.Item($r, 20).Value = $T

    # Per-row values that actually change in this edit.
    $ws.Cells.Item($r, 4).Value = $rd.D
    # Row 17 is brand new, so it needs the same date number format the other
    # "Fecha" cells use (existing rows already have it; re-applying is a no-op).
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 13).Value = $rd.M
    $ws.Cells.Item($r, 14).Value = $rd.N
    $ws.Cells.Item($r, 15).Value = $rd.O
    $ws.Cells.Item($r, 16).Value = $rd.P
    $ws.Cells.Item($r, 17).Value = $rd.Q
    $ws.Cells.Item($r, 19).Value = $rd.S
}
